$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: the word "numberOfAliens" inside
#   numberOfAliens.append(numberOfAliens[i-1] + numberHatched)
# was stored as two runs ("nu" / "mberOfAliens") with a _GoBack bookmark
# sitting between them. Re-finding/replacing that occurrence collapses
# it back into a single "numberOfAliens" run (the visible text itself
# does not change).
# ---------------------------------------------------------------------
$aliensPara = $d.Paragraphs(65).Range

# Skip past the first "numberOfAliens" (the one in "numberOfAliens.append")
# so the Find below lands on the second occurrence - the split one.
$skip = $aliensPara.Duplicate
$skip.Find.Execute("numberOfAliens.append(", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$target = $d.Range($skip.End, $aliensPara.End)
$target.Find.Execute("numberOfAliens", $true, $false, $false, $false, $false, $true, 1, $false, "numberOfAliens", 2)

# ---------------------------------------------------------------------
# Change 2: add a new paragraph after "Big O - O(n)" justifying the
# runtime/Big-O answer for question 2 part 3.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count).Range
$lastPara.InsertParagraphAfter()

$newPara = $d.Paragraphs($d.Paragraphs.Count).Range
$newPara.Text = "Depends on hatching days too as some parts of code run for the length of m and the other part of the if statement runs for n-m. However as n is the most costly, for example as the number of days inputted is the larger number between n and m, as otherwise no eggs would hatch, it thus determines how many times the majority of the code runs. Therefore the Big O value is n. "
